## backwardElimination.xlsx edit
## - Sheet "2" (the last backward-elimination step, Df Model = 2) is removed as a
##   separate tab.
## - Sheet "3" keeps its slot/tab, but its summary cell + row height are refreshed
##   to the content that used to live on sheet "2" (the regression that used the
##   Ones/LangScr/MathScr columns), with the run's Date/Time stamp updated.
## - The other two surviving summary cells (sheets "5" and "4") get their
##   Date/Time stamp refreshed the same way.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$oldDate = "Date:                Wed, 25 Dec 2019"
$newDate = "Date:                Sat, 28 Dec 2019"
$oldTime = "Time:                        23:10:07"
$newTime = "Time:                        20:51:41"

function Update-Timestamp([string]$text) {
    $text = $text.Replace($oldDate, $newDate)
    $text = $text.Replace($oldTime, $newTime)
    return $text
}

$wsFive   = $wb.Worksheets.Item("5")
$wsFour   = $wb.Worksheets.Item("4")
$wsThree  = $wb.Worksheets.Item("3")
$wsTwo    = $wb.Worksheets.Item("2")

# Capture sheet "2"'s summary text + row height before it disappears: this is the
# content that ends up "inherited" by sheet "3" once the dedicated Df-Model=3
# write-up goes away.
$twoSummary   = $wsTwo.Range("B2").Text
$twoRowHeight = $wsTwo.Rows.Item(2).RowHeight

# Refresh the run timestamp baked into each surviving OLS summary blob.
$wsFive.Range("B2").Value  = Update-Timestamp $wsFive.Range("B2").Text
$wsFour.Range("B2").Value  = Update-Timestamp $wsFour.Range("B2").Text

# Sheet "3" takes over sheet "2"'s write-up (with the refreshed timestamp) and
# matches its row height.
$wsThree.Range("B2").Value = Update-Timestamp $twoSummary
$wsThree.Rows.Item(2).RowHeight = $twoRowHeight

# Drop the now-redundant "2" tab.
$wsTwo.Delete()

# Keep the original active tab ("5", the workbook's first sheet) selected.
$wsFive.Activate()
